# Update "want to go" / attendance counts (column F) on the "展览"
# (Exhibitions) and "全部类型" (All types) sheets, regenerated from
# the upstream data source (gh-pages build at 456a3b4).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 16161
$ws.Range("F5").Value = 422
$ws.Range("F6").Value = 13
$ws.Range("F8").Value = 15552
$ws.Range("F9").Value = 65
$ws.Range("F10").Value = 9183
$ws.Range("F11").Value = 449
$ws.Range("F14").Value = 117
$ws.Range("F19").Value = 79
$ws.Range("F24").Value = 1137
$ws.Range("F28").Value = 514
$ws.Range("F30").Value = 43
$ws.Range("F36").Value = 346
$ws.Range("F37").Value = 470
$ws.Range("F39").Value = 5648

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 16161
$ws.Range("F5").Value = 422
$ws.Range("F6").Value = 13
$ws.Range("F8").Value = 15553
$ws.Range("F9").Value = 65
$ws.Range("F10").Value = 9183
$ws.Range("F11").Value = 449
$ws.Range("F14").Value = 117
$ws.Range("F19").Value = 79
$ws.Range("F24").Value = 1137
$ws.Range("F28").Value = 514
$ws.Range("F30").Value = 43
$ws.Range("F38").Value = 346
$ws.Range("F39").Value = 470
$ws.Range("F41").Value = 5648
